$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row above row 50, pushing the existing rows 50-115 down to 51-116.
$ws.Rows(50).Insert()

# Populate the newly inserted row 50 with the new "Cebollín" price record.
$ws.Cells.Item(50, 1).Value = 11
$ws.Cells.Item(50, 2).Value = "Vega Monumental Concepción"
$ws.Cells.Item(50, 3).Value = "Bíobío"
$ws.Cells.Item(50, 4).Value = 45079
$ws.Cells.Item(50, 5).Value = 8
$ws.Cells.Item(50, 6).Value = 100112037
$ws.Cells.Item(50, 7).Value = "Cebollín"
$ws.Cells.Item(50, 8).Value = "Sin especificar"
$ws.Cells.Item(50, 9).Value = "Primera"
$ws.Cells.Item(50, 10).Value = 50
$ws.Cells.Item(50, 11).Value = 3500
$ws.Cells.Item(50, 12).Value = 4000
$ws.Cells.Item(50, 13).Value = 3700
$ws.Cells.Item(50, 14).Value = "$/paquete 36 unidades"
$ws.Cells.Item(50, 15).Value = "Región Metropolitana"
$ws.Cells.Item(50, 16).Value = 103
$ws.Cells.Item(50, 17).Value = 36
$ws.Cells.Item(50, 18).Value = "Hortaliza"
